$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.466.75'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '2.071.63'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.00'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +1.81%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.29'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("E9").Value = '  +2.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0778'
$ws.Range("E10").Value = '  +2.14%  '
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").Value = '2.374.63'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.37'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.68'
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.779'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.20'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = '2.056.50'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Value = '37.352.78'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.61'
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.51'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.46'
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.134'
$ws.Range("E27").Value = '  +5.40%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.82'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  -5.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.10'
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.57'
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  +0.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("E35").Value = '  -2.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.37'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  -4.28%  '
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").Value = '1.491.29'
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0955'
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.02'
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0213'
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.18'
$ws.Range("E46").Value = '  -5.16%  '
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("E48").Value = '  -3.51%  '
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").Value = '2.261.09'
$ws.Range("E51").Value = '  +0.18%  '
